$d = $word.ActiveDocument

# Locate the paragraph that currently ends with
# "b) Yes, each solution will work in all cases" (and holds the _GoBack bookmark).
$targetIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text
    if ($txt -like "*Yes, each solution will work in all cases*") {
        $targetIndex = $i
    }
}

$p38 = $d.Paragraphs.Item($targetIndex)

# Insert two new (plain, non-bold) paragraphs after it: one for item 5's
# heading and one for item 5's answer text.
$p38.Range.InsertParagraphAfter()
$headingPara = $d.Paragraphs.Item($targetIndex + 1)

$headingPara.Range.InsertParagraphAfter()
$answerPara = $d.Paragraphs.Item($targetIndex + 2)

# Fill in the answer paragraph's text first (still unbolded at this point).
$answerIp = $d.Range($answerPara.Range.End - 1, $answerPara.Range.End - 1)
$answerIp.InsertAfter([char]9 + "a) When you choose 12 socks, you have picked enough socks to guarantee " + [char]9 + "you will have at least one matching pair of socks and when you choose 18 " + [char]9 + "socks, you have picked enough socks to guarantee you will have at least one " + [char]9 + "matching pair of socks of each color")

# Fill in the heading paragraph's text and bold it.
$headingIp = $d.Range($headingPara.Range.End - 1, $headingPara.Range.End - 1)
$headingIp.InsertAfter("5) Choose a solution and develop a plan to implement it")
$headingPara.Range.Bold = 1

# Move the _GoBack bookmark so that it now sits between "a) W" and "hen" in the
# newly typed answer paragraph (this also removes it from its old location,
# since bookmark names are unique).
$bookmarkPos = $answerPara.Range.Start + 5
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
